# Rename the sheet to reflect the new "updated on" date.
# Excel automatically re-points the workbook-scoped defined name
# "Socialfaglige_systemer" (which refers to this sheet's A1:E11 range)
# to the new sheet name when the sheet is renamed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Opdateret d. 05-12-2025"
